$wb = $excel.ActiveWorkbook

# --- Update selections on the existing sheets ---
$wsNumero = $wb.Worksheets.Item("Numero spettacoli")
$wsNumero.Activate()
$wsNumero.Range("H5").Select()

$wsIngressi = $wb.Worksheets.Item("Ingressi")
$wsIngressi.Activate()
$wsIngressi.Range("B2:F13").Select()

$wsSpesa = $wb.Worksheets.Item("Spesa del pubblico")
$wsSpesa.Activate()
$wsSpesa.Range("B2:F13").Select()

# --- Add the new "Corretto" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCorretto = $wb.Worksheets.Add($null, $lastSheet)
$wsCorretto.Name = "Corretto"

# --- Header row ---
$wsCorretto.Range("A1").Value = "Date"
$wsCorretto.Range("B1").Value = "Numero spettacoli"
$wsCorretto.Range("C1").Value = "Ingressi"
$wsCorretto.Range("D1").Value = "Spesa del pubblico"

$data = @(
    @(43101, 14048, 2557873, 38990026.829999998),
    @(43132, 14911, 2734855, 37233797.500000007),
    @(43160, 16203, 2973107, 40538564.569999978),
    @(43191, 13106, 2173417, 31807364.410000004),
    @(43221, 10938, 1651762, 28185285.089999992),
    @(43252, 8489, 1293711, 31427737.800000001),
    @(43282, 7336, 1140402, 43536833.089999996),
    @(43313, 5295, 831804, 29351240.070000004),
    @(43344, 5940, 640135, 41239398.109999992),
    @(43374, 10364, 1416863, 42491218.230000012),
    @(43405, 14536, 2524772, 42542731.300000012),
    @(43435, 14640, 2752923, 52612956.360000029),
    @(43466, 13353, 2703477, 47919487.549999997),
    @(43497, 14981, 3053040, 43509811.680000007),
    @(43525, 17129, 3272010, 43349471.219999999),
    @(43556, 10621, 1794711, 25812150.910000011),
    @(43586, 11156, 1747340, 26245363.179999996),
    @(43617, 8444, 1314409, 31951599.989999998),
    @(43647, 6603, 998253, 39053084.380000003),
    @(43678, 4976, 892017, 30347584.98),
    @(43709, 5633, 675616, 45793519.530000009),
    @(43739, 10309, 1500830, 48884917.219999991),
    @(43770, 14859, 2579042, 44371790.359999992),
    @(43800, 14137, 2797637, 57994067.570000015),
    @(43831, 13697, 2756682, 46002045.609999999),
    @(43862, 13225, 2716268, 36430573.399999999),
    @(43891, 477, 80448, 1364026.6199999996),
    @(43922, $null, $null, 74865.259999999995),
    @(43952, $null, $null, 57814.53),
    @(43983, 827, 38099, 1139897.1499999999),
    @(44013, 4834, 314592, 6785651.0000000009),
    @(44044, 4862, 430722, 9519323.2599999998),
    @(44075, 4227, 303784, 7941193.8100000005),
    @(44105, 4378, 262140, 4776139.0399999991),
    @(44136, $null, $null, 37522),
    @(44166, $null, $null, 69800.600000000006),
    @(44197, $null, $null, 77648.84),
    @(44228, $null, $null, 6286.54),
    @(44256, $null, $null, 49324.959999999999),
    @(44287, 60, 3788, 86016.709999999992),
    @(44317, 3461, 233234, 2418386.84),
    @(44348, 7260, 483472, 6709532.669999999),
    @(44378, 10228, 838959, 18030715.34),
    @(44409, 7318, 721574, 17119925.270000003),
    @(44440, 8068, 578038, 16288714.460000001),
    @(44470, 9074, 861013, 19633424.379999999),
    @(44501, 10813, 1222266, 22971482.030000001),
    @(44531, 11812, 1398111, 28891294.100000005),
    @(44562, 6356, 983778, 19269992.599999994),
    @(44593, 7940, 1276450, 20582183.030000001),
    @(44621, 11789, 1860044, 31916586.140000001),
    @(44652, 10480, 1660376, 27439343.019999988),
    @(44682, 11837, 1822010, 30965381.899999995),
    @(44713, 9530, 1411919, 26838089.020000003),
    @(44743, 9531, 1359619, 36162022.699999996),
    @(44774, 6927, 1024471, 26804926.340000004),
    @(44805, 8197, 977016, 29685883.100000005),
    @(44835, 11880, 1613974, 36457300.020000011),
    @(44866, 14007, 2187237, 38729827.100000009),
    @(44896, 15292, 2665646, 52970629.560000002),
)

# --- Data rows ---
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2

    $wsCorretto.Cells.Item($r, 1).Value = $row[0]
    $wsCorretto.Cells.Item($r, 1).NumberFormat = "mmm-yy"

    if ($null -ne $row[1]) { $wsCorretto.Cells.Item($r, 2).Value = $row[1] }
    if ($null -ne $row[2]) { $wsCorretto.Cells.Item($r, 3).Value = $row[2] }
    if ($null -ne $row[3]) { $wsCorretto.Cells.Item($r, 4).Value = $row[3] }
}

# --- Final selection / activation state ---
$wsCorretto.Range("G24").Select()
$wsCorretto.Activate()
